$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format column G for new rows as Text so "04" is stored as a string, not the number 4
$ws.Range("G21:G25").NumberFormat = "@"

$ws.Cells.Item(21, 1).Value = 1100
$ws.Cells.Item(21, 2).Value = 28.88203620910645
$ws.Cells.Item(21, 3).Value = 3.235722541809082
$ws.Cells.Item(21, 4).Value = 11.2804708480835
$ws.Cells.Item(21, 5).Value = 38.84702301025391
$ws.Cells.Item(21, 6).Value = 18476
$ws.Cells.Item(21, 7).Value = "04"

$ws.Cells.Item(22, 1).Value = 1200
$ws.Cells.Item(22, 2).Value = 28.87217330932617
$ws.Cells.Item(22, 3).Value = 3.326181650161743
$ws.Cells.Item(22, 4).Value = 11.16770267486572
$ws.Cells.Item(22, 5).Value = 38.6898307800293
$ws.Cells.Item(22, 6).Value = 18568
$ws.Cells.Item(22, 7).Value = "04"

$ws.Cells.Item(23, 1).Value = 1300
$ws.Cells.Item(23, 2).Value = 28.82461738586426
$ws.Cells.Item(23, 3).Value = 3.360751867294312
$ws.Cells.Item(23, 4).Value = 11.23262977600098
$ws.Cells.Item(23, 5).Value = 39.287841796875
$ws.Cells.Item(23, 6).Value = 18438
$ws.Cells.Item(23, 7).Value = "04"

$ws.Cells.Item(24, 1).Value = 1400
$ws.Cells.Item(24, 2).Value = 28.71969985961914
$ws.Cells.Item(24, 3).Value = 3.431125402450562
$ws.Cells.Item(24, 4).Value = 10.55944156646729
$ws.Cells.Item(24, 5).Value = 39.51337432861328
$ws.Cells.Item(24, 6).Value = 18366
$ws.Cells.Item(24, 7).Value = "04"

$ws.Cells.Item(25, 1).Value = 1500
$ws.Cells.Item(25, 2).Value = 28.67596054077148
$ws.Cells.Item(25, 3).Value = 3.5433030128479
$ws.Cells.Item(25, 4).Value = 10.44325637817383
$ws.Cells.Item(25, 5).Value = 38.67103576660156
$ws.Cells.Item(25, 6).Value = 18392
$ws.Cells.Item(25, 7).Value = "04"
